# Launch and login feature reviewed and fixed
$wb = $excel.ActiveWorkbook

$wsLogin  = $wb.Worksheets.Item("LoginPageContent")
$wsLinked = $wb.Worksheets.Item("LinkedListPageContent")

# --- LoginPageContent: the suite used to cover every scenario twice, once
# submitting the form normally and once via pressing Enter. The Enter-key
# variants were reviewed out, leaving one row per scenario. ---
$wsLogin.Rows.Item(13).Delete()   # valid_login / confirms login using Enter
$wsLogin.Rows.Item(11).Delete()   # Invalid password / confirms login using Enter
$wsLogin.Rows.Item(9).Delete()    # Invalid cred / confirms login using Enter
$wsLogin.Rows.Item(7).Delete()    # Null value in username / presses Enter
$wsLogin.Rows.Item(5).Delete()    # Null value in password / presses Enter
$wsLogin.Rows.Item(3).Delete()    # Null value in cred / presses Enter

# The remaining "Invalid cred" scenario (now row 5) was reworded.
$wsLogin.Range("A5").Value = "Invalid  user"
$wsLogin.Range("E5").Value = "Invalid Username and Password"
$wsLogin.Range("E6").Value = "Invalid Username and Password"

# Page setup was touched (e.g. print preview) on the login sheet.
$wsLogin.PageSetup.Orientation = 1

# --- View/selection state left by the reviewer ---
[void]$wsLinked.Activate()
[void]$wsLinked.Range("C31").Select()

[void]$wsLogin.Activate()
[void]$wsLogin.Range("A7").Select()
